# Generate Report for Handback
# ------------------------------------------------------------------
# This script updates the localization-status workbook to reflect a
# completed handback: the overall status text changes, a "Latest
# Target File" hyperlink + value is recorded for each localized file
# row (zh-cn and de-de sheets), a "Latest Handback File" value is
# filled in, and the "Latest Handback DateTime" is stamped. Column
# widths are widened to fit the new, longer content.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$url43d = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9aaaaf2bdbdea9e53df3ad64c511daf022256b7f/e2e/43d57670-a5ea-4ec0-bcc8-55f43ac70ef1.md"
$url7feb = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9aaaaf2bdbdea9e53df3ad64c511daf022256b7f/e2e/7feb5e31-61bf-419d-9953-b48d1106cc57.md"

$disp43d = "43d57670-a5ea-4ec0-bcc8-55f43ac70ef1.md"
$disp7feb = "7feb5e31-61bf-419d-9953-b48d1106cc57.md"

# ------------------------------------------------------------------
# Overview sheet: refresh the per-language status text + widen the
# status columns (E = zh-cn, F = de-de).
# ------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ------------------------------------------------------------------
# Helper: stamp the "handed back" info onto a language sheet (zh-cn
# or de-de). Columns: A=Source File Name, C=Status, G=Latest Handoff
# File, H=Latest Handoff Datetime, I=Latest Target File, J=Latest
# Handback File, K=Latest Handback DateTime.
# (Positional params only -- this PowerShell host does not bind
# named `-param value` arguments on custom functions.)
# ------------------------------------------------------------------
function Update-LanguageSheet($ws, $handbackFile2, $handbackFile3, $handbackDate2, $handbackDate3) {

    # Status column picks up the new shared text too.
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Latest Target File (I) -- filename + hyperlink back to the source .md
    $ws.Range("I2").Value = $disp43d
    $ws.Range("I2").Font.Name = "Calibri"
    $ws.Range("I2").Font.Size = 11
    $ws.Range("I2").Font.Underline = 2
    $ws.Range("I2").Font.Color = 15570276
    $ws.Hyperlinks.Add($ws.Range("I2"), $url43d, "", "", $disp43d)

    $ws.Range("I3").Value = $disp7feb
    $ws.Range("I3").Font.Name = "Calibri"
    $ws.Range("I3").Font.Size = 11
    $ws.Range("I3").Font.Underline = 2
    $ws.Range("I3").Font.Color = 15570276
    $ws.Hyperlinks.Add($ws.Range("I3"), $url7feb, "", "", $disp7feb)

    # Latest Handback File (J)
    $ws.Range("J2").Value = $handbackFile2
    $ws.Range("J3").Value = $handbackFile3

    # Latest Handback DateTime (K)
    $ws.Range("K2").Value = $handbackDate2
    $ws.Range("K3").Value = $handbackDate3

    # Column widths: C (Status), I (Latest Target File), J (Latest Handback File)
    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
}

# ------------------------------------------------------------------
# zh-cn sheet
# ------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-LanguageSheet $wsZhCn `
    "43d57670-a5ea-4ec0-bcc8-55f43ac70ef1.052e771efdb624257cfb6c357d7fb343d228c1ca.zh-cn.xlf" `
    "7feb5e31-61bf-419d-9953-b48d1106cc57.36a1ff19bb5cd47d87af07ff0c32a72b1e18355f.zh-cn.xlf" `
    "2016-09-02 06:35:54" `
    "2016-09-02 06:35:54"

# ------------------------------------------------------------------
# de-de sheet
# ------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
Update-LanguageSheet $wsDeDe `
    "43d57670-a5ea-4ec0-bcc8-55f43ac70ef1.052e771efdb624257cfb6c357d7fb343d228c1ca.de-de.xlf" `
    "7feb5e31-61bf-419d-9953-b48d1106cc57.36a1ff19bb5cd47d87af07ff0c32a72b1e18355f.de-de.xlf" `
    "2016-09-02 06:36:02" `
    "2016-09-02 06:36:02"
